$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = "Name"

# Row 2
$ws.Cells.Item(2,1).Value = "Jonny"
$ws.Cells.Item(2,2).Value = "Regiment 9"
$ws.Cells.Item(2,3).Value = "1/2015,6/2016"

# Row 3
$ws.Cells.Item(3,2).Value = "GT"
$ws.Cells.Item(3,3).Value = "8/2011,5/2015"

# Row 4
$ws.Cells.Item(4,2).Value = "Regiment 8"
$ws.Cells.Item(4,3).Value = "7/2016,"
$ws.Cells.Item(4,3).NumberFormat = "mmm-yy"

# Row 5
$ws.Cells.Item(5,1).Value = "Name"

# Row 6
$ws.Cells.Item(6,1).Value = "Bobbu"
$ws.Cells.Item(6,2).Value = "Regiment 7"
$ws.Cells.Item(6,3).Value = "9/2005,8/2014"

# Row 7
$ws.Cells.Item(7,2).Value = "UGA"
$ws.Cells.Item(7,3).Value = "8/2000,5/2004"

# Row 8
$ws.Cells.Item(8,2).Value = "Regiment "
$ws.Cells.Item(8,3).Value = "5/2015,9/2017"

# Row 9
$ws.Cells.Item(9,1).Value = "Name"

# Row 10
$ws.Cells.Item(10,1).Value = "Parker"
$ws.Cells.Item(10,2).Value = "Regiment 4"
$ws.Cells.Item(10,3).Value = "3/2015,8/2018"

# Row 11
$ws.Cells.Item(11,2).Value = "Atlanta Symposium"
$ws.Cells.Item(11,3).Value = "5/2019,5/2019"
$ws.Cells.Item(11,3).NumberFormat = "mmm-yy"

# Row 12
$ws.Cells.Item(12,2).Value = "UGA"
$ws.Cells.Item(12,3).Value = "2/2015,5/2016"

# Row 13
$ws.Cells.Item(13,1).Value = "Name"

# Row 14
$ws.Cells.Item(14,1).Value = "Annie"
$ws.Cells.Item(14,2).Value = "Regiment 7"
$ws.Cells.Item(14,3).Value = "6/2004,8/2007"

# Row 15
$ws.Cells.Item(15,2).Value = "UGA"
$ws.Cells.Item(15,3).Value = "8/2000,5/2004"

# Row 16
$ws.Cells.Item(16,2).Value = "Regiment "
$ws.Cells.Item(16,3).Value = "5/2015,"

# Column H width (matches the added <col min="8" max="8" .../> in the diff)
$ws.Columns.Item(8).ColumnWidth = 10.7

# View state: scroll so row 7 is at the top, and select E12 (matches sheetView in the diff)
$ws.Range("E12").Select() | Out-Null
